$d = $word.ActiveDocument

# The original document has a single paragraph:
#   "This is a text for master22"
# followed by the _GoBack bookmark.
#
# Target:
#   Paragraph 1: three runs "S" + "alam" + " master" (=> "Salam master"),
#                still followed by the _GoBack bookmark.
#   Paragraph 2 (new): "This is a text for master22"

$p1 = $d.Paragraphs(1)
$full = $p1.Range

$xml = "<?xml version='1.0'?>" +
       "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
       "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
       "<pkg:xmlData>" +
       "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
       "<w:body>" +
       "<w:p>" +
       "<w:r><w:t>S</w:t></w:r>" +
       "<w:r><w:t>alam</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'> master</w:t></w:r>" +
       "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
       "<w:bookmarkEnd w:id='0'/>" +
       "</w:p>" +
       "</w:body>" +
       "</w:document>" +
       "</pkg:xmlData></pkg:part></pkg:package>"

# Replace paragraph 1's content with the three runs + bookmark above.
# This also splits off the paragraph mark into a new, now-empty,
# second paragraph.
[void]$full.InsertXML($xml)

# Fill the new second paragraph with the original sentence.
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "This is a text for master22"
